$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting (from the last header cell, AC1) onto the
# three new header cells so they pick up the same style index (bold, centered,
# bordered) instead of Excel synthesizing a brand-new style entry.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record (Wins / Losses / Ties) repeated for every player row, 2-42.
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95   # AD
    $ws.Cells.Item($r, 31).Value = 67   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
